$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Fill in the new row's values: newest date plus the same metric values
# as the previous top data row (783.5 / 1112 / 3610)
# Force column A to be treated as plain text so the date string is not
# auto-converted into a date serial value.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-11-27"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# The row Insert operation copies the formatting of the row above (the
# header row), so reset the new data row back to the default/unstyled
# look used by the rest of the data rows.
$ws.Range("A2:D2").Style = "Normal"

$wb.Save()
